$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '24.770.75'
$ws.Range('E2').Value = '  +0.43%  '
$ws.Range('D3').Value = '1.704.36'
$ws.Range('E3').Value = '  +0.31%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.003'
$ws.Range('D4').Style = "Normal"
$ws.Range('E4').Value = '  +0.27%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '316.93'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -0.14%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '1.004'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +0.35%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.3941'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  -0.61%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.4048'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +0.17%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '1.524'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -1.02%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '1.003'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +0.24%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '53.49'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -2.02%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.08888'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +0.71%  '
$ws.Range('E13').Value = '  +0.56%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '23.71'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +1.23%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '8.062'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +5.28%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.00001328'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -0.41%  '
$ws.Range('D17').Value = '1.718.64'
$ws.Range('E17').Value = '  +0.89%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '99.97'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -1.50%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.07048'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -0.79%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '19.75'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -0.41%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '7.079'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +2.54%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '1.002'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +0.17%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '14.48'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +2.06%  '
$ws.Range('D24').Value = '24.760.93'
$ws.Range('E24').Value = '  +0.46%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '3.222'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +3.90%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '2.369'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +1.56%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '22.79'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +1.37%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '162.33'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +1.41%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '8.830'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +15.94%  '
$ws.Range('E30').Value = '  +1.51%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '5.175'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -1.48%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '7.969'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +5.99%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.08921'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +3.72%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.083'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -3.47%  '
$ws.Range('E35').Value = '  +1.80%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '11.09'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -4.68%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.2764'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +0.23%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '14.64'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -1.03%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.02799'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -0.26%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.09188'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +0.85%  '
$ws.Range('B41').Value = 'TrustWalletToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '1.462'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -0.27%  '
$ws.Range('B42').Value = 'TheSandbox'
$ws.Range('C42').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.7732'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -0.57%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '15.99'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +2.52%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.7219'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -0.91%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '2.583'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +1.73%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '4.213'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -0.36%  '
$ws.Range('E47').Value = '  +0.19%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.338'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -2.78%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '140.95'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -0.83%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '90.85'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +1.90%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.07994'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -0.81%  '
